$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("Assets")
$c = $ws3.Range("C2")
$c.Value = "x"
$b = $c.Borders.Item(9)
$b.Color = 14737632
$b.Weight = -4138
Write-Host "done"
